# Weekly price update: insert 5 new rows of "Limón" price data (dated 44448)
# for "Vega Modelo de Temuco", pushing the existing rows 767-781 down to 772-786.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 blank rows starting at row 767 (shifts old rows 767-781 -> 772-786).
$ws.Range("A767:A771").EntireRow.Insert()

# New data rows, columns A..T.
$newRows = @(
    @{ Row = 767; Values = @(10, "Vega Modelo de Temuco", "La Araucanía", 44448, 9, "Fruta", 100102, "Cítricos", 100102003, "Limón", "Sin especificar", "1a amarillo", 400, 7000, 8000, 7500, "$/bandeja 15 kilos", "Región de O'Higgins", 500, 15) },
    @{ Row = 768; Values = @(10, "Vega Modelo de Temuco", "La Araucanía", 44448, 9, "Fruta", 100102, "Cítricos", 100102003, "Limón", "Sin especificar", "1a amarillo", 10, 130000, 130000, 130000, "$/bins (450 kilos)", "Región de O'Higgins", 289, 450) },
    @{ Row = 769; Values = @(10, "Vega Modelo de Temuco", "La Araucanía", 44448, 9, "Fruta", 100102, "Cítricos", 100102003, "Limón", "Sin especificar", "2a amarillo", 90, 6000, 6000, 6000, "$/bandeja 15 kilos", "Región de O'Higgins", 400, 15) },
    @{ Row = 770; Values = @(10, "Vega Modelo de Temuco", "La Araucanía", 44448, 9, "Fruta", 100102, "Cítricos", 100102003, "Limón", "Sin especificar", "2a amarillo", 300, 4000, 4000, 4000, "$/malla 15 kilos", "Región de O'Higgins", 267, 15) },
    @{ Row = 771; Values = @(10, "Vega Modelo de Temuco", "La Araucanía", 44448, 9, "Fruta", 100102, "Cítricos", 100102003, "Limón", "Sutil De Gase", "Primera", 50, 45000, 45000, 45000, "$/caja 24 kilos", "Perú", 1875, 24) }
)

foreach ($entry in $newRows) {
    $r = $entry.Row
    $vals = $entry.Values
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $vals[$i]
    }
}

Write-Output "Inserted 5 new price rows (767-771)."
